$wb = $excel.ActiveWorkbook

# --- Sheet "REGCV1" (sheet10.xml): remove the stray cell N6 = "ß".
#     Row 7 (U7) is untouched and keeps its own row index - this is a
#     content clear, not a row/shift delete. ---
$wsRegcv1 = $wb.Worksheets.Item("REGCV1")
$wsRegcv1.Range("N6").ClearContents()

# --- Sheet "BusROCOF" (sheet2.xml): fix E7 (35 -> 32) and append rows 8-11 ---
$wsBusRocof = $wb.Worksheets.Item("BusROCOF")

# All data in this sheet is stored as text (numbers-as-text), matching the
# existing rows 1-7. Force the text number format before writing so the
# numeric-looking values ("6", "0.002", ...) are preserved verbatim rather
# than being coerced to real numbers.
$wsBusRocof.Range("A8:I11").NumberFormat = "@"
$wsBusRocof.Range("E7").NumberFormat = "@"

$wsBusRocof.Range("E7").Value = "32"

$newRows = @(
    @("6", "1", "BusROCOF_7",  "SG_3", "33", "0.002", "0.02", "60", "0.001"),
    @("7", "1", "BusROCOF_8",  "SG_4", "34", "0.002", "0.02", "60", "0.001"),
    @("8", "1", "BusROCOF_9",  "SG_5", "36", "0.002", "0.02", "60", "0.001"),
    @("9", "1", "BusROCOF_10", "SG_6", "39", "0.002", "0.02", "60", "0.001")
)

$r = 8
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $wsBusRocof.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
